# Update Audience Segment & Contact Category
#
# The "Segment" header/column becomes "Tags", the old "Tags" header/column
# becomes "Category", and the sample data is refreshed to match:
#   - Company "Example Corp" -> "Acme Corp"
#   - old Segment value "Marketing" -> "Lead,VIP" (now under the Tags header)
#   - old Tags value "Lead,VIP" -> "Clients" (now under the Category header)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the header row (E1 was "Segment", F1 was "Tags").
$ws.Range("E1").Value = "Tags"
$ws.Range("F1").Value = "Category"

# Give the relabeled E1 header the same formatting as the other header
# cell that closes out the table (matches F1's right-edge header style).
$ws.Range("F1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

# Refresh the sample row's data for the renamed columns.
$ws.Range("F2").Value = "Clients"
$ws.Range("D2").Value = "Acme Corp"
$ws.Range("E2").Value = "Lead,VIP"

# The blank placeholder cell in the now-unused E column of row 3 is cleared
# entirely (value + formatting), leaving the row one cell shorter.
$ws.Range("E3").Clear() | Out-Null

# Restore the saved selection to D3.
$ws.Range("D3").Select() | Out-Null
